# Refresh the cryptocurrency price/volume snapshot (scheduled scrape).
# $wb / $ws are pre-bound to the open workbook per the harness contract.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.431.08"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.557.21"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'593.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'173.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.555.57"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "'5.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "'27.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "3.004.83"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "67.256.93"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "2.533.85"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "'7.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'11.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "'357.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.24%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'70.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("D28").Value = "2.689.46"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'537.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'8.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.18%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "'157.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'18.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "'2.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.31%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'39.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "'151.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("D49").Value = "'0.568"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  +1.20%  "
